$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add row 31
$ws.Range("A31").Value = 10001
$ws.Range("B31").Value = 110030
$ws.Range("C31").Value = 10030
$ws.Range("D31").Value = "eng"
$ws.Range("E31").Value = $true
$ws.Range("F31").Value = "superadmin"
$ws.Range("G31").Value = "now()"
$ws.Range("H31").Value = "now()"

# Add row 32
$ws.Range("A32").Value = 10001
$ws.Range("B32").Value = 110031
$ws.Range("C32").Value = 10031
$ws.Range("D32").Value = "eng"
$ws.Range("E32").Value = $true
$ws.Range("F32").Value = "superadmin"
$ws.Range("G32").Value = "now()"
$ws.Range("H32").Value = "now()"

# Update selection/view to mimic scroll position in diff
$ws.Range("F30").Select()
$excel.ActiveWindow.ScrollRow = 25
